$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: snapshot donor row formats into a scratch area (rows 101+) ---
# so that in-place reordering below cannot clobber a donor before it is used.
$ws.Range("A1:E1").Copy()
$ws.Range("A101:E101").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A102:E102").PasteSpecial(-4122)
$ws.Range("A10:E10").Copy()
$ws.Range("A103:E103").PasteSpecial(-4122)
$ws.Range("A15:E15").Copy()
$ws.Range("A104:E104").PasteSpecial(-4122)
$ws.Range("A19:E19").Copy()
$ws.Range("A105:E105").PasteSpecial(-4122)
$ws.Range("A22:E22").Copy()
$ws.Range("A106:E106").PasteSpecial(-4122)

# --- Phase 2: clear main area A1:E24 (values + formats) ---
$ws.Range("A1:E24").Clear()

# --- Phase 3: paste formats from scratch rows into final destination rows ---
$ws.Range("A101:E101").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)
$ws.Range("A102:E102").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$ws.Range("A102:E102").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)
$ws.Range("A102:E102").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A102:E102").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A102:E102").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A103:E103").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A104:E104").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A104:E104").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A104:E104").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("A105:E105").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A105:E105").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("A105:E105").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A106:E106").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A106:E106").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

# --- Phase 4: clear scratch area ---
$ws.Range("A101:E106").Clear()

# --- Phase 5: set cell values ---
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "numkey"
$ws.Range("C1").Value = "Gamepad"
$ws.Range("D1").Value = "Alphakey"
$ws.Range("E1").Value = "Function"
$ws.Range("A2").Value = "P1+"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "a"
$ws.Range("E2").Value = "() => document.getElementById(`"p1plus`").click()"
$ws.Range("A3").Value = "P1-"
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = "q"
$ws.Range("E3").Value = "() => document.getElementById(`"p1moins`").click()"
$ws.Range("A4").Value = "P2+"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "k"
$ws.Range("E4").Value = "() => document.getElementById(`"p2plus`").click()"
$ws.Range("A5").Value = "P2-"
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = "o"
$ws.Range("E5").Value = "() => document.getElementById(`"p2moins`").click()"
$ws.Range("A6").Value = "P1P2 zéro"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = "p"
$ws.Range("E6").Value = "resetScores"
$ws.Range("A7").Value = "P1x"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = "z"
$ws.Range("E7").Value = "() => addTime('addButton')"
$ws.Range("A8").Value = "P2x"
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "m"
$ws.Range("E8").Value = "() => addTime('addButton2')"
$ws.Range("A9").Value = "Coup suivant"
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "c"
$ws.Range("E9").Value = "resetToNextValue"
$ws.Range("A10").Value = "Pause Reprendre"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "v"
$ws.Range("E10").Value = "pauseTimer"
$ws.Range("A11").Value = "Nouvelle partie"
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "b"
$ws.Range("E11").Value = "resetTimer"
$ws.Range("A12").Value = "Afficher minuteur"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 13
$ws.Range("D12").Value = "f"
$ws.Range("E12").Value = "toggleVisibility"
$ws.Range("A13").Value = "Réglages"
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = "r"
$ws.Range("E13").Value = "openSettings"
$ws.Range("A14").Value = "Instructions"
$ws.Range("C14").Value = 8
$ws.Range("D14").Value = "i"
$ws.Range("E14").Value = "openInstructions"
$ws.Range("A15").Value = "Afficher boutons"
$ws.Range("C15").Value = 16
$ws.Range("D15").Value = "d"
$ws.Range("E15").Value = "toggleBoutonsRonds"
$ws.Range("A16").Value = "Agrandir"
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = "t"
$ws.Range("E16").Value = "zoomIn"
$ws.Range("A17").Value = "Réduire"
$ws.Range("B17").Value = "5-LP"
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = "y"
$ws.Range("E17").Value = "zoomOut"
$ws.Range("A18").Value = "Plein écran"
$ws.Range("C18").Value = "13-LP"
$ws.Range("D18").Value = "u"
$ws.Range("E18").Value = "toggleFullscreen"
$ws.Range("A19").Value = "Recharger page"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = "w"
$ws.Range("E19").Value = "reloadPage"
$ws.Range("A20").Value = "Passer accueil"
$ws.Range("B20").Value = "9-LP"
$ws.Range("C20").Value = "11-LP"
$ws.Range("D20").Value = "e"
$ws.Range("E20").Value = "hideSplashScreen"

# --- Phase 6: sheet-level properties ---
# Columns C:D already sit at 11.42578125 (inherited from the original B:D group) so
# they are left untouched. Columns A:B need to grow to ~20.7109375 "characters"; the
# ColumnWidth setter here snaps to a coarse pixel grid, so 19.75 is the closest input
# that lands on the nearest achievable width.
$ws.Columns("A:B").ColumnWidth = 19.75
$ws.Range("B20").Select()
